$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.675.57"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "2.200.49"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'258.35"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "'83.62"
$ws.Range("E6").Value = "  +11.51%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D9").Value = "'0.598"
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").Value = "'44.42"
$ws.Range("E10").Value = "  +10.20%  "
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "'7.18"
$ws.Range("E12").Value = "  +5.44%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "2.529.28"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "'14.35"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "2.224.41"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "43.605.54"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "'69.67"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").Value = "'2.36"
$ws.Range("E22").Value = "  +11.04%  "
$ws.Range("D23").Value = "'230.31"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("D24").Value = "'8.99"
$ws.Range("E24").Value = "  -4.23%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "'3.51"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'10.64"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'39.04"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").Value = "'173.56"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").Value = "'0.0856"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0361"
$ws.Range("E37").Value = "  +7.08%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.50"
$ws.Range("E38").Value = "  +6.41%  "
$ws.Range("D39").Value = "'12.53"
$ws.Range("E39").Value = "  +4.90%  "
$ws.Range("D40").Value = "'2.85"
$ws.Range("E40").Value = "  +10.26%  "
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "'62.96"
$ws.Range("E42").Value = "  +6.90%  "
$ws.Range("E43").Value = "  +6.27%  "
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("D45").Value = "'8.37"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "'0.0978"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'99.58"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  +6.32%  "
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").Value = "'0.437"
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("E51").Value = "  +7.63%  "
